$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.381.99'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.73%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.325.29'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.23%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '180.46'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.80%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.620'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.84%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('E9').Value = '  -3.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.64'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.402'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.91%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.904.68'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('E13').Value = '  -1.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '66.444.53'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.70'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.91%  '
$ws.Range('E16').Value = '  -1.97%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.335.89'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '436.04'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.61%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.66'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.52'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.78%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.58'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.79%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.17'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.29%  '
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.518'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.75%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000116'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.81%  '
$ws.Range('E26').Value = '  +1.29%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.04'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.39%  '
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('E29').Value = '  -2.86%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.74'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.00%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.23'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.42%  '
$ws.Range('E33').Value = '  -2.89%  '
$ws.Range('E34').Value = '  -4.46%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '160.42'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.89%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.47'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '27.69'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.79'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.82%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.829.56'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.78%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.792'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.77%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.42'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.34%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.15'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.86%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.22'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.96%  '
$ws.Range('E44').Value = '  -3.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '24.08'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.08%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.33'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.50%  '
$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '324.40'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.27%  '
$ws.Range('E48').Value = '  -4.13%  '
$ws.Range('E49').Value = '  +1.07%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.14'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.48%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.970'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.01%  '
